$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.533.06"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.411.71"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +8.46%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +13.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.74"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.656"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.03%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.667"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +11.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.33"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0945"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.64"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.31"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +16.40%  "

$ws.Range("E15").Value = "  +2.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.775.62"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +8.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.421.11"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +8.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.565.18"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.50"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.43%  "

$ws.Range("E20").Value = "  +5.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.18"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.51"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "272.32"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +17.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.46"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.64"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +7.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.02"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.51%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.97"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "23.06"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +10.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "177.91"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.68%  "

$ws.Range("E31").Value = "  +0.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.88"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.76%  "

$ws.Range("E33").Value = "  +4.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0939"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.97"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.57%  "

$ws.Range("E36").Value = "  +6.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.88"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.34%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.08"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.76%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0372"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.109"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +20.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.62"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +24.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.76"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +26.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.236"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.26"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.78"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.45%  "

$ws.Range("E47").Value = "  +0.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.73"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +15.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.70"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +7.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "87.72"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +67.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.32"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.09%  "
